# "docs: first batch of tb report"
# Fill in hours + work-done description for the two rows that were still
# missing them at the bottom of the "Journal de travail" table, then move
# the on-screen selection down to the new last row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Journal de travail")

# Row 80 (2023-07-18, Implémentation): add the missing hours and extend
# the work-done comment with the new items that were completed.
$ws.Range("C80").Value = 5
$ws.Range("D80").Value = "Correction de bugs, doc swagger, feature d'envoie d'event à plusieurs modules"

# Row 81 (2023-07-18, Rédaction): add hours + the work-done comment.
$ws.Range("C81").Value = 4
$ws.Range("D81").Value = "Structure du rapport, modules et validation"

# Reflect the new bottom-of-table position in the sheet view/selection.
$ws.Activate()
$ws.Range("D91").Select()
